# Auto-generated Excel COM-interop script to apply scheduled-runner market-data updates
# to the Garuda_Profits workbook (per-sheet "Sheets" tables: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 33283
$ws.Range("J10").Value = 33283
$ws.Range("L10").Value = 33283
$ws.Range("N10").Value = -33869
$ws.Range("H13").Value = 61145.43
$ws.Range("J13").Value = 70003.60000000001
$ws.Range("L13").Value = 70003.60000000001
$ws.Range("N13").Value = -70341.60000000001
$ws.Range("H33").Value = 186.125
$ws.Range("I33").Value = 193.83018
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 193.83018
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = 35.16981999999999
$ws.Range("N33").Value = -508
$ws.Range("H64").Value = 3019.3333
$ws.Range("I64").Value = 2540
$ws.Range("J64").Value = 3259
$ws.Range("K64").Value = 2540
$ws.Range("L64").Value = 3259
$ws.Range("M64").Value = -2292
$ws.Range("N64").Value = -3755
$ws.Range("H67").Value = 3019.3333
$ws.Range("I67").Value = 2540
$ws.Range("J67").Value = 3259
$ws.Range("K67").Value = 2540
$ws.Range("L67").Value = 3259
$ws.Range("M67").Value = -1682
$ws.Range("N67").Value = -4975
$ws.Range("H87").Value = 31000
$ws.Range("I87").Value = 30000
$ws.Range("J87").Value = 31666.666
$ws.Range("K87").Value = 30000
$ws.Range("L87").Value = 31666.666
$ws.Range("M87").Value = -28752
$ws.Range("N87").Value = -34162.666
$ws.Range("H90").Value = 31000
$ws.Range("I90").Value = 30000
$ws.Range("J90").Value = 31666.666
$ws.Range("K90").Value = 90000
$ws.Range("L90").Value = 94999.99800000001
$ws.Range("M90").Value = -83760
$ws.Range("N90").Value = -107479.998
$ws.Range("H137").Value = 2093.3845
$ws.Range("I137").Value = 1380.4
$ws.Range("J137").Value = 2539
$ws.Range("K137").Value = 4141.200000000001
$ws.Range("L137").Value = 7617
$ws.Range("M137").Value = -1591.200000000001
$ws.Range("N137").Value = -12717

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2127.2415
$ws.Range("I61").Value = 1604
$ws.Range("J61").Value = 3121.4
$ws.Range("K61").Value = 1604
$ws.Range("L61").Value = 3121.4
$ws.Range("M61").Value = -1392
$ws.Range("N61").Value = -3545.4
$ws.Range("H74").Value = 1198.8334
$ws.Range("I74").Value = 1007.11536
$ws.Range("J74").Value = 2445
$ws.Range("K74").Value = 1007.11536
$ws.Range("L74").Value = 2445
$ws.Range("M74").Value = -133.11536
$ws.Range("N74").Value = -4193
$ws.Range("H77").Value = 1198.8334
$ws.Range("I77").Value = 1007.11536
$ws.Range("J77").Value = 2445
$ws.Range("K77").Value = 5035.5768
$ws.Range("L77").Value = 12225
$ws.Range("M77").Value = -667.5767999999998
$ws.Range("N77").Value = -20961
$ws.Range("H132").Value = 6880.7393
$ws.Range("I132").Value = 7547.6113
$ws.Range("J132").Value = 4480
$ws.Range("K132").Value = 22642.8339
$ws.Range("L132").Value = 13440
$ws.Range("M132").Value = -20112.8339
$ws.Range("N132").Value = -18500
$ws.Range("H136").Value = 2127.2415
$ws.Range("I136").Value = 1604
$ws.Range("J136").Value = 3121.4
$ws.Range("K136").Value = 4812
$ws.Range("L136").Value = 9364.200000000001
$ws.Range("M136").Value = -2262
$ws.Range("N136").Value = -14464.2

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 204
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 204
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 204
$ws.Range("H134").Value = 19341.5
$ws.Range("I134").Value = 25681.928
$ws.Range("J134").Value = 2011
$ws.Range("K134").Value = 77045.784
$ws.Range("L134").Value = 6033
$ws.Range("M134").Value = -74510.784
$ws.Range("N134").Value = -11103
$ws.Range("N7").Value = -430
$ws.Range("M7").ClearContents()

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2866.5454
$ws.Range("I62").Value = 2288.5715
$ws.Range("J62").Value = 3878
$ws.Range("K62").Value = 2288.5715
$ws.Range("L62").Value = 3878
$ws.Range("M62").Value = -1664.5715
$ws.Range("N62").Value = -5126
$ws.Range("H65").Value = 2866.5454
$ws.Range("I65").Value = 2288.5715
$ws.Range("J65").Value = 3878
$ws.Range("K65").Value = 11442.8575
$ws.Range("L65").Value = 19390
$ws.Range("M65").Value = -8322.8575
$ws.Range("N65").Value = -25630
$ws.Range("H132").Value = 1966.5778
$ws.Range("I132").Value = 1872.5883
$ws.Range("J132").Value = 2257.0908
$ws.Range("K132").Value = 5617.7649
$ws.Range("L132").Value = 6771.2724
$ws.Range("M132").Value = -3087.7649
$ws.Range("N132").Value = -11831.2724

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 900
$ws.Range("I5").Value = 900
$ws.Range("K5").Value = 900
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10338
$ws.Range("H12").Value = 3835850.2
$ws.Range("I12").Value = 4744467
$ws.Range("J12").Value = 1110000
$ws.Range("K12").Value = 4744467
$ws.Range("L12").Value = 1110000
$ws.Range("M12").Value = -4744327
$ws.Range("N12").Value = -1110280
$ws.Range("H126").Value = 2552.625
$ws.Range("I126").Value = 2778.3635
$ws.Range("K126").Value = 8335.0905
$ws.Range("M126").Value = -5865.0905
$ws.Range("H132").Value = 32564.576
$ws.Range("I132").Value = 36501.07
$ws.Range("K132").Value = 109503.21
$ws.Range("M132").Value = -106973.21
$ws.Range("M5").Value = -788

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1632.5
$ws.Range("J68").Value = 2700
$ws.Range("L68").Value = 2700
$ws.Range("N68").Value = -4198
$ws.Range("H71").Value = 1632.5
$ws.Range("J71").Value = 2700
$ws.Range("L71").Value = 13500
$ws.Range("N71").Value = -20988
$ws.Range("H132").Value = 6936.816
$ws.Range("I132").Value = 10204.565
$ws.Range("J132").Value = 1926.2667
$ws.Range("K132").Value = 30613.695
$ws.Range("L132").Value = 5778.800099999999
$ws.Range("M132").Value = -28083.695
$ws.Range("N132").Value = -10838.8001
$ws.Range("H133").Value = 18883.766
$ws.Range("J133").Value = 18883.766
$ws.Range("L133").Value = 18883.766
$ws.Range("N133").Value = -23943.766

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 9119
$ws.Range("J69").Value = 9119
$ws.Range("L69").Value = 9119
$ws.Range("H70").Value = 26848.666
$ws.Range("J70").Value = 26848.666
$ws.Range("L70").Value = 26848.666
$ws.Range("N70").Value = -27478.666
$ws.Range("H72").Value = 9119
$ws.Range("J72").Value = 9119
$ws.Range("L72").Value = 27357
$ws.Range("H73").Value = 26848.666
$ws.Range("J73").Value = 26848.666
$ws.Range("L73").Value = 26848.666
$ws.Range("N73").Value = -29032.666
$ws.Range("H132").Value = 1096.25
$ws.Range("I132").Value = 1090.1621
$ws.Range("J132").Value = 1128.4286
$ws.Range("K132").Value = 3270.4863
$ws.Range("L132").Value = 3385.2858
$ws.Range("M132").Value = -740.4863
$ws.Range("N132").Value = -8445.2858
$ws.Range("H136").Value = 4966.4546
$ws.Range("I136").Value = 6181.1665
$ws.Range("J136").Value = 1727.2222
$ws.Range("K136").Value = 18543.4995
$ws.Range("L136").Value = 5181.6666
$ws.Range("M136").Value = -15993.4995
$ws.Range("N136").Value = -10281.6666
$ws.Range("N69").Value = -10617
$ws.Range("N72").Value = -34845
